# Ajout de l'activité du jour
# Adds a new journal entry (row 59) to the "Feuil1" worksheet, mirroring
# the formatting of the previous entry (row 58).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 58
$newRow = 59

# Copy the formatting (styles, wrap text, date number format, row layout)
# from the last existing entry so the new row reuses the same cell styles
# instead of creating new ones.
$ws.Range("A" + $lastRow + ":C" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New entry's date: 2018-03-28
$ws.Range("A" + $newRow).Value = 43187

# Description of the day's work.
$ws.Range("B" + $newRow).Value = "J'ai réglé deux problèmes que j'avais. Quand j'ajoutais un article dans mon panier d'une taille différente que la première seléctionnée il m'ajoutais bien mon article dans mon panier mais il m'enlevait une unité à la première taille de l'article que j'ai ajouté dans le panier. J'ai changé une variable qui reprennait l'id de mon article. L'autre problème c'était quand je cliquais sur mon bouton pour ajouter un article qui était plus en stock il me faisait une erreur, j'ai ajouté une condition en plus pour contrôler ça et afficher un message qui disait qu'il n'y avait plus la taille dispo pour cette article."

# Duration of the work session.
$ws.Range("C" + $newRow).Value = "3 périodes"

# Grow the row so the wrapped text is fully visible (matches the taller
# row height Excel would have applied automatically for this entry).
$ws.Rows.Item($newRow).RowHeight = 90

# Move the selection to the next empty row, as left behind by the author.
[void]$ws.Range("C60").Select()
